$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new column before column N ---
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()

# New column N takes the width of column M (its left neighbour) to match
# the format Excel copies on column insert.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Update the selection on this sheet and make it the active tab.
$ws.Activate()
$ws.Range("R10").Select()

$wb.Save()
